$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Advance the order date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Fix bug: prices were double what they should be (exceeded the requested
# amount pulled from google drive) - halve them.
$ws.Range("D29").Value = 56.2
$ws.Range("D30").Value = 93.59999999999999
